$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 149, shifting existing rows 149:294 down to 150:295
$ws.Rows("149:149").Insert()

# Populate the newly inserted row 149 with the new data record
$ws.Range("A149").Value = 4
$ws.Range("B149").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C149").Value = "Los Lagos"
$ws.Range("D149").Value = 44781
$ws.Range("E149").Value = 10
$ws.Range("F149").Value = 100112017
$ws.Range("G149").Value = "Apio"
$ws.Range("H149").Value = "Americana (o)"
$ws.Range("I149").Value = "Primera"
$ws.Range("J149").Value = 20
$ws.Range("K149").Value = 14000
$ws.Range("L149").Value = 14000
$ws.Range("M149").Value = 14000
$ws.Range("N149").Value = "`$/docena de matas"
$ws.Range("O149").Value = "Región de Coquimbo"
$ws.Range("P149").Value = 2333
$ws.Range("Q149").Value = 6
$ws.Range("R149").Value = "Hortaliza"
